$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price (D) and 1h volume-change (E) columns with the
# latest scraped snapshot. Some "Price" strings are plain decimal
# numbers (e.g. "580.51") which Excel would otherwise auto-convert to
# a Number when assigned via .Value; force those specific cells to
# stay Text (matching the source data which is always a literal
# string), then restore the Normal style so no formatting residue
# is left behind.

$ws.Range("D2").Value = '66.525.22'
$ws.Range("E2").Value = '  -1.12%  '
$ws.Range("D3").Value = '3.451.16'
$ws.Range("E3").Value = '  -0.86%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '580.51'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.28%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '175.89'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.60%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.599'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.16%  '
$ws.Range("D9").Value = '3.451.31'
$ws.Range("E9").Value = '  -0.93%  '
$ws.Range("E10").Value = '  -2.34%  '
$ws.Range("E11").Value = '  -3.22%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.418'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.71%  '
$ws.Range("D13").Value = '4.044.70'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '30.50'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.51%  '
$ws.Range("E15").Value = '  -3.54%  '
$ws.Range("D16").Value = '66.521.59'
$ws.Range("E16").Value = '  -1.19%  '
$ws.Range("D18").Value = '3.451.37'
$ws.Range("E18").Value = '  -0.87%  '
$ws.Range("E19").Value = '  -4.02%  '
$ws.Range("E20").Value = '  -3.15%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '375.89'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.34%  '
$ws.Range("E22").Value = '  -2.74%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.998'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.09%  '
$ws.Range("E24").Value = '  -0.15%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '71.15'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.23%  '
$ws.Range("E26").Value = '  -1.72%  '
$ws.Range("E27").Value = '  -2.98%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.80'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.60%  '
$ws.Range("E29").Value = '  -2.00%  '
$ws.Range("E30").Value = '  -0.03%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.84'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.10%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '23.95'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.79%  '
$ws.Range("E33").Value = '  -3.75%  '
$ws.Range("E34").Value = '  -6.03%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.08%  '
$ws.Range("E36").Value = '  -4.47%  '
$ws.Range("E37").Value = '  -4.98%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '159.43'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.86%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.876'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.57%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '27.37'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.18%  '
$ws.Range("E41").Value = '  -5.18%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.62'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.67%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.49'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.36%  '
$ws.Range("E44").Value = '  -3.94%  '
$ws.Range("D45").Value = '2.686.51'
$ws.Range("E45").Value = '  -5.55%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0694'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.73%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '25.21'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.18%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '40.17'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.31%  '
$ws.Range("E49").Value = '  -1.93%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '320.21'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.84%  '
$ws.Range("E51").Value = '  -4.04%  '
